$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.Value = "'" + $val
    $r.Style = "Normal"
}

Set-TextCell $ws "D2" "45.369.62"
Set-TextCell $ws "E2" "  +2.57%  "
Set-TextCell $ws "D3" "2.424.30"
Set-TextCell $ws "E3" "  -0.34%  "
Set-TextCell $ws "E4" "  +0.02%  "
Set-TextCell $ws "D5" "318.85"
Set-TextCell $ws "E5" "  +3.67%  "
Set-TextCell $ws "D6" "103.14"
Set-TextCell $ws "E6" "  +3.73%  "
Set-TextCell $ws "E7" "  +0.49%  "
Set-TextCell $ws "E8" "  -0.03%  "
Set-TextCell $ws "E9" "  +6.04%  "
Set-TextCell $ws "D10" "35.58"
Set-TextCell $ws "E10" "  +1.08%  "
Set-TextCell $ws "D11" "0.0804"
Set-TextCell $ws "E11" "  +0.49%  "
Set-TextCell $ws "E12" "  -2.12%  "
Set-TextCell $ws "D13" "18.22"
Set-TextCell $ws "E13" "  -2.65%  "
Set-TextCell $ws "D14" "7.08"
Set-TextCell $ws "E14" "  +2.07%  "
Set-TextCell $ws "D15" "2.805.28"
Set-TextCell $ws "E15" "  -0.05%  "
Set-TextCell $ws "D16" "2.414.26"
Set-TextCell $ws "E16" "  +0.14%  "
Set-TextCell $ws "E17" "  +1.21%  "
Set-TextCell $ws "D18" "45.298.73"
Set-TextCell $ws "E18" "  +2.55%  "
Set-TextCell $ws "D19" "12.21"
Set-TextCell $ws "E19" "  -0.96%  "
Set-TextCell $ws "D21" "0.0₃0924"
Set-TextCell $ws "E21" "  +1.94%  "
Set-TextCell $ws "D22" "68.95"
Set-TextCell $ws "E22" "  +0.46%  "
Set-TextCell $ws "D23" "244.64"
Set-TextCell $ws "E23" "  +1.77%  "
Set-TextCell $ws "E24" "  -0.99%  "
Set-TextCell $ws "E25" "  +0.84%  "
Set-TextCell $ws "D27" "25.76"
Set-TextCell $ws "E27" "  +1.82%  "
Set-TextCell $ws "B28" "Toncoin"
Set-TextCell $ws "C28" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell $ws "D28" "2.19"
Set-TextCell $ws "E28" "  -1.29%  "
Set-TextCell $ws "B29" "Cosmos"
Set-TextCell $ws "C29" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell $ws "D29" "9.60"
Set-TextCell $ws "E29" "  +0.95%  "
Set-TextCell $ws "D30" "49.35"
Set-TextCell $ws "E30" "  +2.57%  "
Set-TextCell $ws "D31" "32.98"
Set-TextCell $ws "E31" "  +0.22%  "
Set-TextCell $ws "D32" "20.35"
Set-TextCell $ws "E32" "  +9.21%  "
Set-TextCell $ws "E33" "  +3.84%  "
Set-TextCell $ws "D34" "5.21"
Set-TextCell $ws "E34" "  +0.78%  "
Set-TextCell $ws "E35" "  +0.25%  "
Set-TextCell $ws "D36" "0.0767"
Set-TextCell $ws "E36" "  -0.05%  "
Set-TextCell $ws "B37" "RenderToken"
Set-TextCell $ws "C37" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws "D37" "4.46"
Set-TextCell $ws "E37" "  -2.04%  "
Set-TextCell $ws "B38" "ARBITRUM"
Set-TextCell $ws "C38" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell $ws "D38" "1.87"
Set-TextCell $ws "E38" "  -3.05%  "
Set-TextCell $ws "D39" "2.87"
Set-TextCell $ws "E39" "  -2.13%  "
Set-TextCell $ws "D40" "125.43"
Set-TextCell $ws "E40" "  -4.12%  "
Set-TextCell $ws "B41" "Stellar"
Set-TextCell $ws "C41" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell $ws "D41" "0.110"
Set-TextCell $ws "E41" "  +0.53%  "
Set-TextCell $ws "B42" "WEMIXToken"
Set-TextCell $ws "C42" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D42" "2.21"
Set-TextCell $ws "E42" "  -3.45%  "
Set-TextCell $ws "D43" "20.59"
Set-TextCell $ws "E43" "  -2.01%  "
Set-TextCell $ws "E44" "  +1.07%  "
Set-TextCell $ws "D45" "1.932.14"
Set-TextCell $ws "E45" "  -1.31%  "
Set-TextCell $ws "E46" "  -2.94%  "
Set-TextCell $ws "E47" "  +1.46%  "
Set-TextCell $ws "E48" "  -2.30%  "
Set-TextCell $ws "D49" "1.79"
Set-TextCell $ws "E49" "  +8.64%  "
Set-TextCell $ws "D50" "76.79"
Set-TextCell $ws "E50" "  +4.38%  "
Set-TextCell $ws "D51" "53.79"
Set-TextCell $ws "E51" "  +0.52%  "
